$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "se actualiza el diccionario de datos" - update the data dictionary:
# drop the Creditos / Contra_creditos columns (M:N) and the
# Porcentaje_Ejecucion column (originally T, now shifted to R) from the
# sheet. Deleting the entire columns shifts everything after them left,
# which updates the header row, the shared-string-backed data row and
# the dimension/spans automatically.
$ws.Range("M1:N1").EntireColumn.Delete() | Out-Null
$ws.Range("R1").EntireColumn.Delete() | Out-Null

# Restore the view the author left the workbook in: zoomed to 130% with
# Q4 selected.
$excel.ActiveWindow.Zoom = 130
$ws.Range("Q4").Select() | Out-Null
